# Rename the embedded picture's docPr/cNvPr "name" attribute on the four
# Pearson/BTec logo inline pictures that live in the document's headers and
# footers:
#   - header (first page)  : id=1  image1.jpg -> image2.jpg
#   - header (default)     : id=3  image1.jpg -> image2.jpg
#   - footer (first page)  : id=2  image2.png -> image1.png
#   - footer (default)     : id=4  image2.png -> image1.png
#
# InlineShape objects don't expose a high-level property that reaches the
# nested <pic:cNvPr name="..."/> element, so we round-trip the shape's own
# Range through WordOpenXML: read it, rewrite the `name="..."` attribute
# (it occurs exactly twice - once on wp:docPr, once on pic:cNvPr - and both
# need the same new value), and feed it back with InsertXML.
#
# Re-importing a header/footer's XML like this makes Word materialize a
# "Header"/"Footer" paragraph style (with its linked character style) in
# styles.xml if one isn't already defined, mirroring what headers/footers
# implicitly use. Since that's not part of the intended edit, remove any
# such style(s) that get newly added right after each InsertXML call so
# the style list ends up exactly as it started.

function Rename-ShapeName($shape, [string]$oldName, [string]$newName) {
    $d = $word.ActiveDocument
    $beforeCount = $d.Styles.Count
    $beforeNames = @()
    for ($i = 1; $i -le $beforeCount; $i++) {
        $beforeNames += $d.Styles.Item($i).NameLocal
    }

    $range = $shape.Range
    $xml = $range.WordOpenXML
    $newXml = $xml.Replace('name="' + $oldName + '"', 'name="' + $newName + '"')
    $range.InsertXML($newXml)

    $afterCount = $d.Styles.Count
    if ($afterCount -gt $beforeCount) {
        for ($i = $afterCount; $i -ge 1; $i--) {
            $st = $d.Styles.Item($i)
            if (-not ($beforeNames -contains $st.NameLocal)) {
                $st.Delete()
            }
        }
    }
}

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# Headers: Item(1) is the default header, Item(2) is the first-page header.
$headerDefault = $section.Headers.Item(1)
$headerFirst = $section.Headers.Item(2)

Rename-ShapeName $headerDefault.Range.InlineShapes.Item(1) "image1.jpg" "image2.jpg"
Rename-ShapeName $headerFirst.Range.InlineShapes.Item(1) "image1.jpg" "image2.jpg"

# Footers: Item(1) is the default footer, Item(2) is the first-page footer.
$footerDefault = $section.Footers.Item(1)
$footerFirst = $section.Footers.Item(2)

Rename-ShapeName $footerDefault.Range.InlineShapes.Item(1) "image2.png" "image1.png"
Rename-ShapeName $footerFirst.Range.InlineShapes.Item(1) "image2.png" "image1.png"
